$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dataset values (after adding falling/struggle/walkingToRunning classes
# and re-running 10-fold cross validation). Rows 2-21 (20 data rows) replace
# the previous rows 2-22 (21 data rows).
$data = @(
    @(14.36146354675293, -2.813630223274229, 2.31547799706459),
    @(7.119595527648926, -10.00965690612793, 2.071655035018921),
    @(-0.3092890381813076, -7.561092853546142, 3.192377448081971),
    @(-0.1974980831146222, -4.94487106800079, 4.721822917461397),
    @(-2.194161117076871, -4.593446969985962, 3.157028853893282),
    @(-0.8586132526397767, -10.09085631370543, -3.983258485794054),
    @(2.208070576190949, -11.74357312917708, -5.92425370216368),
    @(0.8976666927337558, -6.909027695655827, 1.685739278793339),
    @(-3.808973312377933, -9.732912957668304, 4.339319378137589),
    @(1.211441993713379, -4.946407794952393, 2.633680582046509),
    @(-1.650709629058838, -4.344925820827484, 4.874097138643265),
    @(-3.371967196464541, -3.051161766052245, 5.898134231567384),
    @(-2.243164837360378, -0.9728509187698318, 7.228664159774784),
    @(1.634144544601426, -3.087035417556748, 5.197757840156568),
    @(0.4866030216217228, -5.888972640037538, 5.150876790285098),
    @(-3.873410999774934, 2.929124236106873, 4.643965721130372),
    @(3.477634161710739, 24.22354996204378, -8.259647905826574),
    @(3.516963958740234, -0.7771883010864258, 4.514358997344971),
    @(-0.7894209027290375, 0.4619972705841073, 8.847737967967991),
    @(-0.04974877834320068, -2.565077662467957, 6.396960973739624)
)

$rowCount = $data.Count
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Remove the now-obsolete trailing row (previously row 22) since the dataset shrank by one row
$lastOldRow = $rowCount + 2
$ws.Range("A$lastOldRow`:C$lastOldRow").Delete()
